$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.167270660400391
$ws.Range("B1").Value = 2.437687635421753
$ws.Range("D1").Value = 2.368618249893188
$ws.Range("E1").Value = 1.233832001686096
